$wb = $excel.ActiveWorkbook

# Sheet: ALC (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("H98").Value = 622
$ws.Range("I98").Value = 637.1923
$ws.Range("K98").Value = 637.1923
$ws.Range("M98").Value = 860.8077
$ws.Range("H100").Value = 791.6667
$ws.Range("I100").Value = 783.9286
$ws.Range("K100").Value = 783.9286
$ws.Range("M100").Value = -242.9286
$ws.Range("H101").Value = 3293.6667
$ws.Range("J101").Value = 4665
$ws.Range("L101").Value = 13995
$ws.Range("N101").Value = -17239
$ws.Range("H103").Value = 560.75
$ws.Range("I103").Value = 482
$ws.Range("J103").Value = 797
$ws.Range("K103").Value = 1446
$ws.Range("L103").Value = 2391
$ws.Range("M103").Value = -860
$ws.Range("N103").Value = -3563
$ws.Range("H113").Value = 3548.3103
$ws.Range("I113").Value = 3099.0667
$ws.Range("J113").Value = 4029.6428
$ws.Range("K113").Value = 3099.0667
$ws.Range("L113").Value = 4029.6428
$ws.Range("M113").Value = 154.9333000000001
$ws.Range("N113").Value = -10537.6428
$ws.Range("H122").Value = 622
$ws.Range("I122").Value = 637.1923
$ws.Range("K122").Value = 1911.5769
$ws.Range("M122").Value = 538.4231
$ws.Range("H132").Value = 18682.172
$ws.Range("I132").Value = 5848.7144
$ws.Range("J132").Value = 21890.535
$ws.Range("K132").Value = 17546.1432
$ws.Range("L132").Value = 65671.605
$ws.Range("M132").Value = -15016.1432
$ws.Range("N132").Value = -70731.605
$ws.Range("H135").Value = 2515.5
$ws.Range("I135").Value = 898.52
$ws.Range("K135").Value = 8086.68
$ws.Range("M135").Value = -5551.68
$ws.Range("H138").Value = 3620.75
$ws.Range("I138").Value = 1574.6666
$ws.Range("J138").Value = 4069.8901
$ws.Range("K138").Value = 4723.9998
$ws.Range("L138").Value = 12209.6703
$ws.Range("M138").Value = 416.0002000000004
$ws.Range("N138").Value = -22489.6703
# Sheet: ARM (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 7733.6294
$ws.Range("I61").Value = 8841.450000000001
$ws.Range("J61").Value = 4568.4287
$ws.Range("K61").Value = 8841.450000000001
$ws.Range("L61").Value = 4568.4287
$ws.Range("M61").Value = -8629.450000000001
$ws.Range("N61").Value = -4992.4287
$ws.Range("H74").Value = 1476.2858
$ws.Range("I74").Value = 1000.3571
$ws.Range("J74").Value = 2428.1428
$ws.Range("K74").Value = 1000.3571
$ws.Range("L74").Value = 2428.1428
$ws.Range("M74").Value = -126.3570999999999
$ws.Range("N74").Value = -4176.1428
$ws.Range("H77").Value = 1476.2858
$ws.Range("I77").Value = 1000.3571
$ws.Range("J77").Value = 2428.1428
$ws.Range("K77").Value = 5001.7855
$ws.Range("L77").Value = 12140.714
$ws.Range("M77").Value = -633.7855
$ws.Range("N77").Value = -20876.714
$ws.Range("H102").Value = 473971.7
$ws.Range("I102").Value = 528199.5600000001
$ws.Range("J102").Value = 3996.6667
$ws.Range("K102").Value = 528199.5600000001
$ws.Range("L102").Value = 3996.6667
$ws.Range("M102").Value = -526577.5600000001
$ws.Range("N102").Value = -7240.6667
$ws.Range("H122").Value = 4656.2
$ws.Range("I122").Value = 2038.6154
$ws.Range("J122").Value = 9517.429
$ws.Range("K122").Value = 6115.8462
$ws.Range("L122").Value = 28552.287
$ws.Range("M122").Value = -3665.8462
$ws.Range("N122").Value = -33452.287
$ws.Range("H132").Value = 16451.705
$ws.Range("I132").Value = 29056.285
$ws.Range("K132").Value = 87168.855
$ws.Range("M132").Value = -84638.855
$ws.Range("H136").Value = 7733.6294
$ws.Range("I136").Value = 8841.450000000001
$ws.Range("J136").Value = 4568.4287
$ws.Range("K136").Value = 26524.35
$ws.Range("L136").Value = 13705.2861
$ws.Range("M136").Value = -23974.35
$ws.Range("N136").Value = -18805.2861
$ws.Range("H138").Value = 64714.5
$ws.Range("J138").Value = 64714.5
$ws.Range("L138").Value = 64714.5
$ws.Range("N138").Value = -74994.5
# Sheet: BSM (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("H61").Value = 47500
$ws.Range("J61").Value = 50000
$ws.Range("L61").Value = 50000
$ws.Range("N61").Value = -50626
$ws.Range("H105").Value = 2033.1875
$ws.Range("I105").Value = 1810.1538
$ws.Range("K105").Value = 1810.1538
$ws.Range("M105").Value = -63.15380000000005
# Sheet: CRP (index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 6342.1836
$ws.Range("I31").Value = 2633.2222
$ws.Range("K31").Value = 2633.2222
$ws.Range("M31").Value = -2338.2222
$ws.Range("H34").Value = 6342.1836
$ws.Range("I34").Value = 2633.2222
$ws.Range("K34").Value = 2633.2222
$ws.Range("M34").Value = -2431.2222
$ws.Range("H36").Value = 25907
$ws.Range("J36").Value = 31666
$ws.Range("L36").Value = 31666
$ws.Range("N36").Value = -32442
$ws.Range("H39").Value = 7625
$ws.Range("I39").Value = 7625
$ws.Range("K39").Value = 7625
$ws.Range("M39").Value = -7234
$ws.Range("H40").Value = 25907
$ws.Range("J40").Value = 31666
$ws.Range("L40").Value = 31666
$ws.Range("N40").Value = -31986
$ws.Range("H49").Value = 7625
$ws.Range("I49").Value = 7625
$ws.Range("K49").Value = 7625
$ws.Range("M49").Value = -7443
$ws.Range("H99").Value = 7518.2383
$ws.Range("J99").Value = 7804.9414
$ws.Range("L99").Value = 7804.9414
$ws.Range("N99").Value = -10800.9414
$ws.Range("H126").Value = 7518.2383
$ws.Range("J126").Value = 7804.9414
$ws.Range("L126").Value = 23414.8242
$ws.Range("N126").Value = -28354.8242
$ws.Range("H141").Value = 82753.89
$ws.Range("J141").Value = 84074.94
$ws.Range("L141").Value = 84074.94
$ws.Range("N141").Value = -94434.94
# Sheet: CUL (index 5)
$ws = $wb.Worksheets.Item(5)
$ws.Range("H98").Value = 1368.625
$ws.Range("J98").Value = 1658.3334
$ws.Range("L98").Value = 4975.0002
$ws.Range("N98").Value = -7971.0002
$ws.Range("H126").Value = 17010
$ws.Range("I126").Value = 1030
$ws.Range("J126").Value = 25000
$ws.Range("K126").Value = 3090
$ws.Range("L126").Value = 75000
$ws.Range("M126").Value = 1850
$ws.Range("N126").Value = -84880
$ws.Range("H129").Value = 1120.9166
$ws.Range("I129").Value = 917.44446
$ws.Range("J129").Value = 1731.3334
$ws.Range("K129").Value = 2752.33338
$ws.Range("L129").Value = 5194.0002
$ws.Range("M129").Value = 2247.66662
$ws.Range("N129").Value = -15194.0002
$ws.Range("H131").Value = 13266381
$ws.Range("J131").Value = 15161442
$ws.Range("L131").Value = 45484326
$ws.Range("N131").Value = -45494406
# Sheet: GSM (index 6)
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 5687955.5
$ws.Range("I70").Value = 7581274
$ws.Range("K70").Value = 7581274
$ws.Range("M70").Value = -7581004
$ws.Range("H73").Value = 5687955.5
$ws.Range("I73").Value = 7581274
$ws.Range("K73").Value = 7581274
$ws.Range("M73").Value = -7580338
$ws.Range("H122").Value = 651324.2
$ws.Range("I122").Value = 1102254
$ws.Range("K122").Value = 3306762
$ws.Range("M122").Value = -3304312
$ws.Range("H132").Value = 576726.1
$ws.Range("I132").Value = 203635.6
$ws.Range("K132").Value = 610906.8
$ws.Range("M132").Value = -608376.8
# Sheet: LTW (index 7)
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 3473.7827
$ws.Range("I16").Value = 1972.1875
$ws.Range("K16").Value = 1972.1875
$ws.Range("M16").Value = -1802.1875
$ws.Range("H40").Value = 3139.1667
$ws.Range("I40").Value = 3306.25
$ws.Range("J40").Value = 1802.5
$ws.Range("K40").Value = 3306.25
$ws.Range("L40").Value = 1802.5
$ws.Range("M40").Value = -3170.25
$ws.Range("N40").Value = -2074.5
$ws.Range("H51").Value = 37624.75
$ws.Range("J51").Value = 49500
$ws.Range("L51").Value = 49500
$ws.Range("N51").Value = -50456
$ws.Range("H61").Value = 3224.5
$ws.Range("I61").Value = 966
$ws.Range("K61").Value = 966
$ws.Range("M61").Value = -764
$ws.Range("H113").Value = 3224.5
$ws.Range("I113").Value = 966
$ws.Range("K113").Value = 966
$ws.Range("M113").Value = 1204
$ws.Range("H132").Value = 4958.3335
$ws.Range("I132").Value = 5234.875
$ws.Range("J132").Value = 4405.25
$ws.Range("K132").Value = 15704.625
$ws.Range("L132").Value = 13215.75
$ws.Range("M132").Value = -13174.625
$ws.Range("N132").Value = -18275.75
$ws.Range("H136").Value = 2743.7026
$ws.Range("I136").Value = 1938.1786
$ws.Range("J136").Value = 5249.778
$ws.Range("K136").Value = 5814.5358
$ws.Range("L136").Value = 15749.334
$ws.Range("M136").Value = -3264.5358
$ws.Range("N136").Value = -20849.334
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
# Sheet: WVR (index 8)
$ws = $wb.Worksheets.Item(8)
$ws.Range("H41").Value = 17311
$ws.Range("J41").Value = 17311
$ws.Range("L41").Value = 17311
$ws.Range("N41").Value = -18091
$ws.Range("H96").Value = 7764.1665
$ws.Range("I96").Value = 6987
$ws.Range("K96").Value = 6987
$ws.Range("M96").Value = -5614
$ws.Range("H107").Value = 5379.1816
$ws.Range("I107").Value = 6268.8
$ws.Range("J107").Value = 3472.8572
$ws.Range("K107").Value = 18806.4
$ws.Range("L107").Value = 10418.5716
$ws.Range("M107").Value = -16886.4
$ws.Range("N107").Value = -14258.5716
$ws.Range("H122").Value = 3363.4348
$ws.Range("J122").Value = 1641.5
$ws.Range("L122").Value = 4924.5
$ws.Range("N122").Value = -9824.5
$ws.Range("H132").Value = 16667646
$ws.Range("I132").Value = 1218.8334
$ws.Range("J132").Value = 41667290
$ws.Range("K132").Value = 3656.5002
$ws.Range("L132").Value = 125001870
$ws.Range("M132").Value = -1126.5002
$ws.Range("N132").Value = -125006930
